$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # quality_comparison
$ws2 = $wb.Worksheets.Item(2)   # computational_comparison

# --- Sheet1 (quality_comparison): give C1 a top+bottom border, D1 a right+top+bottom border ---
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item(8).LineStyle = 1   # xlEdgeTop      (top only     == pre-existing border 2)
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight     (top+right    == pre-existing border 3)
$d1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom    (top+right+bottom == pre-existing border 5)

# --- Sheet2 (computational_comparison): same border treatment for C1/D1 and F1/G1 ---
# Re-use the already-created styles by copying formats, so no extra style entries are minted.
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)   # xlPasteFormats

$c1.Copy()
$ws2.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

$d1.Copy()
$ws2.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Text changes: "fedcore" -> "approach" ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Numeric sign-of-zero cleanups: -0 -> 0 ---
$ws1.Range("D4").Value = 0
$ws1.Range("D5").Value = 0
$ws1.Range("D12").Value = 0

# --- Remove the stray empty inline-string cell G5 on sheet2 ---
$ws2.Range("G5").ClearContents()
